$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that differ between row 22 and row 23 and need to be swapped
$cols = @("A", "B", "D", "E", "F", "G", "H", "I", "Q", "R")

foreach ($col in $cols) {
    $cell22 = $ws.Range($col + "22")
    $cell23 = $ws.Range($col + "23")

    $val22 = $cell22.Value2
    $val23 = $cell23.Value2

    $cell22.Value2 = $val23
    $cell23.Value2 = $val22
}
